# Auto-generated script applying scheduled market-data refresh to Kujata_Profits workbook
# Updates currentAveragePrice* / Leve*Price* / Leve*Profit* columns (H-N) for affected leve rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 478.625
$ws.Range("I8").Value = 478.625
$ws.Range("K8").Value = 1435.875
$ws.Range("M8").Value = -1296.875
# Row 64
$ws.Range("H64").Value = 4297.3335
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 4000
$ws.Range("N64").Value = -4496
# Row 67
$ws.Range("H67").Value = 4297.3335
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 4000
$ws.Range("N67").Value = -5716
# Row 74
$ws.Range("H74").Value = 4996.4
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# Row 76
$ws.Range("H76").Value = 5766.6665
$ws.Range("I76").Value = 6000
$ws.Range("K76").Value = 6000
$ws.Range("M76").Value = -5685
# Row 77
$ws.Range("H77").Value = 4996.4
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
# Row 79
$ws.Range("H79").Value = 5766.6665
$ws.Range("I79").Value = 6000
$ws.Range("K79").Value = 6000
$ws.Range("M79").Value = -4908
# Row 88
$ws.Range("H88").Value = 950642.7
$ws.Range("J88").Value = 1544144.6
$ws.Range("L88").Value = 1544144.6
$ws.Range("N88").Value = -1544956.6
# Row 91
$ws.Range("H91").Value = 950642.7
$ws.Range("J91").Value = 1544144.6
$ws.Range("L91").Value = 1544144.6
$ws.Range("N91").Value = -1546952.6
# Row 113
$ws.Range("H113").Value = 13335729
$ws.Range("I113").Value = 16669001
$ws.Range("J113").Value = 2638
$ws.Range("K113").Value = 16669001
$ws.Range("L113").Value = 2638
$ws.Range("M113").Value = -16665747
$ws.Range("N113").Value = -9146
# Row 138
$ws.Range("H138").Value = 2690.1648
$ws.Range("I138").Value = 2297.1667
$ws.Range("J138").Value = 2749.8608
$ws.Range("K138").Value = 6891.500100000001
$ws.Range("L138").Value = 8249.582399999999
$ws.Range("M138").Value = -1751.500100000001
$ws.Range("N138").Value = -18529.5824

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11061.512
$ws.Range("I32").Value = 8196.313
$ws.Range("J32").Value = 22353.766
$ws.Range("K32").Value = 8196.313
$ws.Range("L32").Value = 22353.766
$ws.Range("M32").Value = -7909.313
$ws.Range("N32").Value = -22927.766
# Row 132
$ws.Range("H132").Value = 2359.8276
$ws.Range("I132").Value = 1754.9756
$ws.Range("J132").Value = 3818.5881
$ws.Range("K132").Value = 5264.9268
$ws.Range("L132").Value = 11455.7643
$ws.Range("M132").Value = -2734.9268
$ws.Range("N132").Value = -16515.7643

$ws = $wb.Worksheets.Item("BSM")
# Row 75
$ws.Range("H75").Value = 9017.549999999999
$ws.Range("J75").Value = 12724.25
$ws.Range("L75").Value = 12724.25
$ws.Range("N75").Value = -14596.25
# Row 78
$ws.Range("H78").Value = 9017.549999999999
$ws.Range("J78").Value = 12724.25
$ws.Range("L78").Value = 38172.75
$ws.Range("N78").Value = -47532.75
# Row 107
$ws.Range("H107").Value = 893.4211
$ws.Range("I107").Value = 854.0625
$ws.Range("K107").Value = 854.0625
$ws.Range("M107").Value = 1065.9375

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 5066.4814
$ws.Range("I58").Value = 835.5714
$ws.Range("K58").Value = 835.5714
$ws.Range("M58").Value = -632.5714
# Row 62
$ws.Range("H62").Value = 12502220
$ws.Range("I62").Value = 2351.5
$ws.Range("J62").Value = 33335334
$ws.Range("K62").Value = 2351.5
$ws.Range("L62").Value = 33335334
$ws.Range("M62").Value = -1727.5
$ws.Range("N62").Value = -33336582
# Row 65
$ws.Range("H65").Value = 12502220
$ws.Range("I65").Value = 2351.5
$ws.Range("J65").Value = 33335334
$ws.Range("K65").Value = 11757.5
$ws.Range("L65").Value = 166676670
$ws.Range("M65").Value = -8637.5
$ws.Range("N65").Value = -166682910
# Row 74
$ws.Range("H74").Value = 30000
$ws.Range("I74").Value = 20000
$ws.Range("J74").Value = 33333.332
$ws.Range("K74").Value = 20000
$ws.Range("L74").Value = 33333.332
$ws.Range("M74").Value = -19126
$ws.Range("N74").Value = -35081.332
# Row 77
$ws.Range("H77").Value = 30000
$ws.Range("I77").Value = 20000
$ws.Range("J77").Value = 33333.332
$ws.Range("K77").Value = 60000
$ws.Range("L77").Value = 99999.99600000001
$ws.Range("M77").Value = -55632
$ws.Range("N77").Value = -108735.996
# Row 102
$ws.Range("H102").Value = 25000
$ws.Range("J102").Value = 25000
$ws.Range("L102").Value = 25000
$ws.Range("N102").Value = -29868
# Row 108
$ws.Range("H108").Value = 34208
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").ClearContents()
# Row 136
$ws.Range("H136").Value = 5066.4814
$ws.Range("I136").Value = 835.5714
$ws.Range("K136").Value = 2506.7142
$ws.Range("M136").Value = 43.28579999999965

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 41
$ws.Range("I2").Value = 41.8
$ws.Range("J2").Value = 39.666668
$ws.Range("K2").Value = 250.8
$ws.Range("L2").Value = 238.000008
$ws.Range("M2").Value = -137.8
$ws.Range("N2").Value = -464.000008
# Row 12
$ws.Range("H12").Value = 177.375
$ws.Range("I12").Value = 262.5
$ws.Range("K12").Value = 787.5
$ws.Range("M12").Value = -614.5
# Row 17
$ws.Range("H17").Value = 246.5
$ws.Range("I17").Value = 195
$ws.Range("J17").Value = 272.25
$ws.Range("K17").Value = 585
$ws.Range("L17").Value = 816.75
$ws.Range("M17").Value = -416
$ws.Range("N17").Value = -1154.75
# Row 34
$ws.Range("H34").Value = 4763727.5
$ws.Range("J34").Value = 7695037
$ws.Range("L34").Value = 23085111
$ws.Range("N34").Value = -23085279
# Row 39
$ws.Range("H39").Value = 4159.0586
$ws.Range("J39").Value = 4206.933
$ws.Range("L39").Value = 12620.799
$ws.Range("N39").Value = -13208.799
# Row 55
$ws.Range("H55").Value = 1624.75
$ws.Range("I55").Value = 749.5
$ws.Range("J55").Value = 2500
$ws.Range("K55").Value = 2248.5
$ws.Range("L55").Value = 7500
$ws.Range("M55").Value = -2071.5
$ws.Range("N55").Value = -7854
# Row 117
$ws.Range("H117").Value = 744.6923
$ws.Range("J117").Value = 909.2
$ws.Range("L117").Value = 2727.6
$ws.Range("N117").Value = -9611.6
# Row 129
$ws.Range("H129").Value = 19842364
$ws.Range("I129").Value = 41667124
$ws.Range("J129").Value = 6411742
$ws.Range("K129").Value = 125001372
$ws.Range("L129").Value = 19235226
$ws.Range("M129").Value = -124996372
$ws.Range("N129").Value = -19245226

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3436.2727
$ws.Range("J80").Value = 4371.2856
$ws.Range("L80").Value = 4371.2856
$ws.Range("N80").Value = -6367.2856
# Row 83
$ws.Range("H83").Value = 3436.2727
$ws.Range("J83").Value = 4371.2856
$ws.Range("L83").Value = 21856.428
$ws.Range("N83").Value = -31840.428
# Row 107
$ws.Range("H107").Value = 3846744.5
$ws.Range("I107").Value = 6411007.5
$ws.Range("J107").Value = 350
$ws.Range("K107").Value = 6411007.5
$ws.Range("L107").Value = 350
$ws.Range("M107").Value = -6409087.5
$ws.Range("N107").Value = -4190
# Row 122
$ws.Range("H122").Value = 1338
$ws.Range("I122").Value = 1007
$ws.Range("K122").Value = 3021
$ws.Range("M122").Value = -571
# Row 124
$ws.Range("H124").Value = 56000
$ws.Range("J124").Value = 56000
$ws.Range("L124").Value = 56000
$ws.Range("N124").Value = -65820

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 967.2273
$ws.Range("I16").Value = 967.2273
$ws.Range("K16").Value = 967.2273
$ws.Range("M16").Value = -797.2273

$ws = $wb.Worksheets.Item("WVR")
# Row 125
$ws.Range("H125").Value = 28759.8
$ws.Range("J125").Value = 28759.8
$ws.Range("L125").Value = 28759.8
$ws.Range("N125").Value = -38599.8
